$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The June block grew from 20 days to 22 days, so 2 new rows must be
# inserted right after the existing June rows (old row 21, the last
# June entry) to make room before the May block that currently starts
# at row 22.
$ws.Rows.Item(22).Resize(2).Insert()

$data = @"
2,1,5140.01,6,2025,06/2025
3,2,27797.8,6,2025,06/2025
4,3,25236.75,6,2025,06/2025
5,4,36189.46,6,2025,06/2025
6,5,21643.35,6,2025,06/2025
7,6,15533.91,6,2025,06/2025
8,7,19035.4,6,2025,06/2025
9,8,4942.9,6,2025,06/2025
10,9,23856.87,6,2025,06/2025
11,10,31200,6,2025,06/2025
12,11,33164.36,6,2025,06/2025
13,12,38843.63,6,2025,06/2025
14,13,22435.29,6,2025,06/2025
15,14,27509.4,6,2025,06/2025
16,15,4621.42,6,2025,06/2025
17,16,28514.4,6,2025,06/2025
18,17,18098.79,6,2025,06/2025
19,18,27037.48,6,2025,06/2025
20,19,17122.67,6,2025,06/2025
21,20,17622.85,6,2025,06/2025
22,21,13501.01,6,2025,06/2025
23,22,4054.5,6,2025,06/2025
24,1,17056.87,5,2025,05/2025
25,2,22786.63,5,2025,05/2025
26,3,29255,5,2025,05/2025
27,4,4370,5,2025,05/2025
28,5,35402.64,5,2025,05/2025
29,6,28248.43,5,2025,05/2025
30,7,31437.91,5,2025,05/2025
31,8,27732.15,5,2025,05/2025
32,9,25508.17,5,2025,05/2025
33,10,14802.01,5,2025,05/2025
34,11,9776.9,5,2025,05/2025
35,12,31420.44,5,2025,05/2025
36,13,26531.86,5,2025,05/2025
37,14,36574.18,5,2025,05/2025
38,15,33940.79,5,2025,05/2025
39,16,30498.76,5,2025,05/2025
40,17,15120.8,5,2025,05/2025
41,18,8085.01,5,2025,05/2025
42,19,26923.67,5,2025,05/2025
43,20,33090.3,5,2025,05/2025
44,21,25994.16,5,2025,05/2025
45,22,27147.29,5,2025,05/2025
46,23,23820.29,5,2025,05/2025
47,24,34922.72,5,2025,05/2025
48,25,3824.9,5,2025,05/2025
49,26,27831.77,5,2025,05/2025
50,27,24156.31,5,2025,05/2025
51,28,21718.63,5,2025,05/2025
52,29,19445.75,5,2025,05/2025
53,30,26951.12,5,2025,05/2025
54,31,19327.9,5,2025,05/2025
55,1,19371.37,4,2025,04/2025
56,2,29218.71,4,2025,04/2025
57,3,19819.3,4,2025,04/2025
58,4,25399.91,4,2025,04/2025
59,5,18509.5,4,2025,04/2025
60,6,5823.21,4,2025,04/2025
61,7,21764.19,4,2025,04/2025
62,8,38515.34,4,2025,04/2025
63,9,21733.46,4,2025,04/2025
64,10,25715.4,4,2025,04/2025
65,11,22833.25,4,2025,04/2025
66,12,21596.15,4,2025,04/2025
67,13,6451,4,2025,04/2025
68,14,24085.6,4,2025,04/2025
69,15,15644.59,4,2025,04/2025
70,16,18372.2,4,2025,04/2025
71,17,23446.19,4,2025,04/2025
72,18,15474.35,4,2025,04/2025
73,19,9164.799999999999,4,2025,04/2025
74,20,2719,4,2025,04/2025
75,21,15666.86,4,2025,04/2025
76,22,15219.82,4,2025,04/2025
77,23,19895.36,4,2025,04/2025
78,24,19632.82,4,2025,04/2025
79,25,19212.07,4,2025,04/2025
80,26,14442,4,2025,04/2025
81,27,5289.3,4,2025,04/2025
82,28,25680.18,4,2025,04/2025
83,29,21976.8,4,2025,04/2025
84,30,24773.88,4,2025,04/2025
85,1,13803.81,3,2025,03/2025
86,2,4705.9,3,2025,03/2025
87,3,16680.45,3,2025,03/2025
88,4,7657.8,3,2025,03/2025
89,5,37984.79,3,2025,03/2025
90,6,19791.3,3,2025,03/2025
91,7,51663.12,3,2025,03/2025
92,8,14212.1,3,2025,03/2025
93,9,8836.049999999999,3,2025,03/2025
94,10,38138.8,3,2025,03/2025
95,11,30745.57,3,2025,03/2025
96,12,44416.44,3,2025,03/2025
97,13,20045.79,3,2025,03/2025
98,14,23677.64,3,2025,03/2025
99,15,23636.62,3,2025,03/2025
100,16,4642,3,2025,03/2025
101,17,37009.74,3,2025,03/2025
102,18,15851.63,3,2025,03/2025
103,19,33313.16,3,2025,03/2025
104,20,23394.13,3,2025,03/2025
105,21,17765.24,3,2025,03/2025
106,22,21750.9,3,2025,03/2025
107,23,5923.91,3,2025,03/2025
108,24,23370.46,3,2025,03/2025
109,25,27227.58,3,2025,03/2025
110,26,20282.52,3,2025,03/2025
111,27,19695.38,3,2025,03/2025
112,28,17250.36,3,2025,03/2025
113,29,15973.05,3,2025,03/2025
114,30,6906.9,3,2025,03/2025
115,31,24267.47,3,2025,03/2025
"@

$lines = $data -split "`r?`n" | Where-Object { $_.Trim() -ne "" }

foreach ($line in $lines) {
    $parts = $line -split ","
    $row = [int]$parts[0]
    $dia = [double]$parts[1]
    $total = [double]$parts[2]
    $mes = [double]$parts[3]
    $ano = [double]$parts[4]
    $periodo = $parts[5]

    $ws.Cells.Item($row, 1).Value = $dia
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $mes
    $ws.Cells.Item($row, 4).Value = $ano
    $ws.Cells.Item($row, 5).Value = $periodo
}

Write-Output "rows updated: $($lines.Count)"
